$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "26.303.36"
$ws.Range("D3").Value = "1.590.31"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D5").Value = "212.21"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.0609"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "19.35"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.814.67"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.610.52"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "64.40"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "26.313.26"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("D20").Value = "211.89"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "9.00"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").Value = "145.18"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "7.04"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "15.18"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "0.0500"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "1.338.91"
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").Value = "0.601"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "1.05"
$ws.Range("E39").Value = "  -15.69%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "61.89"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "1.726.57"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").Value = "0.0980"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").Value = "  -0.32%  "
